$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Sun Jan 14 16:23:33 EST 2024"
$ws.Range("B3").Value = "Sun Jan 14 16:23:45 EST 2024"
$ws.Range("B4").Value = "Sun Jan 14 16:23:57 EST 2024"
$ws.Range("B5").Value = "Sun Jan 14 16:24:09 EST 2024"
$ws.Range("B6").Value = "Sun Jan 14 16:24:21 EST 2024"
$ws.Range("B7").Value = "Sun Jan 14 16:24:33 EST 2024"
$ws.Range("B8").Value = "Sun Jan 14 16:24:45 EST 2024"
$ws.Range("B9").Value = "Sun Jan 14 16:24:57 EST 2024"
$ws.Range("B10").Value = "Sun Jan 14 16:25:10 EST 2024"
$ws.Range("B11").Value = "Sun Jan 14 16:25:22 EST 2024"
$ws.Range("B12").Value = "Sun Jan 14 16:25:33 EST 2024"
$ws.Range("B13").Value = "Sun Jan 14 16:25:45 EST 2024"
$ws.Range("B14").Value = "Sun Jan 14 16:25:57 EST 2024"
$ws.Range("B15").Value = "Sun Jan 14 16:26:09 EST 2024"
$ws.Range("B16").Value = "Sun Jan 14 16:26:21 EST 2024"
$ws.Range("B17").Value = "Sun Jan 14 16:26:33 EST 2024"
$ws.Range("B18").Value = "Sun Jan 14 16:26:45 EST 2024"
$ws.Range("B19").Value = "Sun Jan 14 16:26:57 EST 2024"
$ws.Range("B20").Value = "Sun Jan 14 16:27:09 EST 2024"
$ws.Range("B28").Value = "Sun Jan 14 16:27:22 EST 2024"
$ws.Range("B29").Value = "Sun Jan 14 16:27:34 EST 2024"
$ws.Range("B30").Value = "Sun Jan 14 16:27:46 EST 2024"
$ws.Range("B31").Value = "Sun Jan 14 16:27:58 EST 2024"
$ws.Range("B32").Value = "Sun Jan 14 16:28:10 EST 2024"
$ws.Range("B33").Value = "Sun Jan 14 16:28:22 EST 2024"
$ws.Range("A34").Value = "Pass"
$ws.Range("A34").Style = "Normal"
$ws.Range("B34").Value = "Sun Jan 14 16:28:33 EST 2024"
$ws.Range("B34").Style = "Normal"
$ws.Range("A35").Value = "Pass"
$ws.Range("A35").Style = "Normal"
$ws.Range("B35").Value = "Sun Jan 14 16:28:46 EST 2024"
$ws.Range("B35").Style = "Normal"
$ws.Range("A36").Value = "Pass"
$ws.Range("A36").Style = "Normal"
$ws.Range("B36").Value = "Sun Jan 14 16:28:58 EST 2024"
$ws.Range("B36").Style = "Normal"
$ws.Range("A37").Value = "Pass"
$ws.Range("A37").Style = "Normal"
$ws.Range("B37").Value = "Sun Jan 14 16:29:10 EST 2024"
$ws.Range("B37").Style = "Normal"
$ws.Range("A38").Value = "Pass"
$ws.Range("A38").Style = "Normal"
$ws.Range("B38").Value = "Sun Jan 14 16:29:22 EST 2024"
$ws.Range("B38").Style = "Normal"
$ws.Range("A39").Value = "Pass"
$ws.Range("A39").Style = "Normal"
$ws.Range("B39").Value = "Sun Jan 14 16:29:34 EST 2024"
$ws.Range("B39").Style = "Normal"
$ws.Range("A40").Value = "Pass"
$ws.Range("A40").Style = "Normal"
$ws.Range("B40").Value = "Sun Jan 14 16:29:46 EST 2024"
$ws.Range("B40").Style = "Normal"
$ws.Range("A41").Value = "Pass"
$ws.Range("A41").Style = "Normal"
$ws.Range("B41").Value = "Sun Jan 14 16:29:58 EST 2024"
$ws.Range("B41").Style = "Normal"
$ws.Range("A42").Value = "Pass"
$ws.Range("A42").Style = "Normal"
$ws.Range("B42").Value = "Sun Jan 14 16:30:10 EST 2024"
$ws.Range("B42").Style = "Normal"
$ws.Range("A43").Value = "Pass"
$ws.Range("A43").Style = "Normal"
$ws.Range("B43").Value = "Sun Jan 14 16:30:22 EST 2024"
$ws.Range("B43").Style = "Normal"
$ws.Range("A44").Value = "Pass"
$ws.Range("A44").Style = "Normal"
$ws.Range("B44").Value = "Sun Jan 14 16:30:34 EST 2024"
$ws.Range("B44").Style = "Normal"
$ws.Range("A45").Value = "Pass"
$ws.Range("A45").Style = "Normal"
$ws.Range("B45").Value = "Sun Jan 14 16:30:46 EST 2024"
$ws.Range("B45").Style = "Normal"
$ws.Range("A46").Value = "Pass"
$ws.Range("A46").Style = "Normal"
$ws.Range("B46").Value = "Sun Jan 14 16:30:58 EST 2024"
$ws.Range("B46").Style = "Normal"
$ws.Range("A47").Value = "Pass"
$ws.Range("A47").Style = "Normal"
$ws.Range("B47").Value = "Sun Jan 14 16:31:10 EST 2024"
$ws.Range("B47").Style = "Normal"
$ws.Range("A48").Value = "Pass"
$ws.Range("A48").Style = "Normal"
$ws.Range("B48").Value = "Sun Jan 14 16:31:22 EST 2024"
$ws.Range("B48").Style = "Normal"
$ws.Range("A49").Value = "Pass"
$ws.Range("A49").Style = "Normal"
$ws.Range("B49").Value = "Sun Jan 14 16:31:34 EST 2024"
$ws.Range("B49").Style = "Normal"
$ws.Range("A50").Value = "Pass"
$ws.Range("A50").Style = "Normal"
$ws.Range("B50").Value = "Sun Jan 14 16:31:46 EST 2024"
$ws.Range("B50").Style = "Normal"
$ws.Range("A51").Value = "Pass"
$ws.Range("A51").Style = "Normal"
$ws.Range("B51").Value = "Sun Jan 14 16:31:58 EST 2024"
$ws.Range("B51").Style = "Normal"
$ws.Range("A52").Value = "Pass"
$ws.Range("A52").Style = "Normal"
$ws.Range("B52").Value = "Sun Jan 14 16:32:10 EST 2024"
$ws.Range("B52").Style = "Normal"
$ws.Range("A53").Value = "Pass"
$ws.Range("A53").Style = "Normal"
$ws.Range("B53").Value = "Sun Jan 14 16:32:22 EST 2024"
$ws.Range("B53").Style = "Normal"
$ws.Range("A54").Value = "Pass"
$ws.Range("A54").Style = "Normal"
$ws.Range("B54").Value = "Sun Jan 14 16:32:34 EST 2024"
$ws.Range("B54").Style = "Normal"
